# Auto-generated script to apply the Odin_Profits market-data refresh edit
# Updates columns H-N (current market price / profit calculations) across rows
# in all 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# Sheet ALC (index 1), row 46
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(46, 8).Value = 0  # H46: 3000 -> 0
$ws.Cells.Item(46, 10).Value = 0  # J46: 3000 -> 0
$ws.Cells.Item(46, 12).Value = 0  # L46: 9000 -> 0
$ws.Cells.Item(46, 14).ClearContents()  # N46: -9238 -> (removed)

# Sheet ALC (index 1), row 60
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(60, 8).Value = 0  # H60: 3000 -> 0
$ws.Cells.Item(60, 10).Value = 0  # J60: 3000 -> 0
$ws.Cells.Item(60, 12).Value = 0  # L60: 9000 -> 0
$ws.Cells.Item(60, 14).ClearContents()  # N60: -9968 -> (removed)

# Sheet ALC (index 1), row 80
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(80, 8).Value = 1745.2941  # H80: 1357.9445 -> 1745.2941
$ws.Cells.Item(80, 9).Value = 1090.909  # I80: 857.5454999999999 -> 1090.909
$ws.Cells.Item(80, 10).Value = 2945  # J80: 2144.2856 -> 2945
$ws.Cells.Item(80, 11).Value = 3272.727  # K80: 2572.6365 -> 3272.727
$ws.Cells.Item(80, 12).Value = 8835  # L80: 6432.8568 -> 8835
$ws.Cells.Item(80, 13).Value = -2274.727  # M80: -1574.6365 -> -2274.727
$ws.Cells.Item(80, 14).Value = -10831  # N80: -8428.856800000001 -> -10831

# Sheet ALC (index 1), row 83
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(83, 8).Value = 1745.2941  # H83: 1357.9445 -> 1745.2941
$ws.Cells.Item(83, 9).Value = 1090.909  # I83: 857.5454999999999 -> 1090.909
$ws.Cells.Item(83, 10).Value = 2945  # J83: 2144.2856 -> 2945
$ws.Cells.Item(83, 11).Value = 9818.181  # K83: 7717.9095 -> 9818.181
$ws.Cells.Item(83, 12).Value = 26505  # L83: 19298.5704 -> 26505
$ws.Cells.Item(83, 13).Value = -4826.181  # M83: -2725.9095 -> -4826.181
$ws.Cells.Item(83, 14).Value = -36489  # N83: -29282.5704 -> -36489

# Sheet ALC (index 1), row 137
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(137, 8).Value = 7251.231  # H137: 5262.4736 -> 7251.231
$ws.Cells.Item(137, 9).Value = 13652.75  # I137: 6033.2 -> 13652.75
$ws.Cells.Item(137, 11).Value = 40958.25  # K137: 18099.6 -> 40958.25
$ws.Cells.Item(137, 13).Value = -38408.25  # M137: -15549.6 -> -38408.25

# Sheet ALC (index 1), row 141
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(141, 8).Value = 4600  # H141: 4170 -> 4600
$ws.Cells.Item(141, 9).Value = 2500  # I141: 2475 -> 2500
$ws.Cells.Item(141, 11).Value = 7500  # K141: 7425 -> 7500
$ws.Cells.Item(141, 13).Value = -2320  # M141: -2245 -> -2320

# Sheet ARM (index 2), row 32
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 1526354.9  # H32: 1614797.1 -> 1526354.9
$ws.Cells.Item(32, 9).Value = 3568.7463  # I32: 3701.5625 -> 3568.7463
$ws.Cells.Item(32, 10).Value = 18530800  # J32: 22236820 -> 18530800
$ws.Cells.Item(32, 11).Value = 3568.7463  # K32: 3701.5625 -> 3568.7463
$ws.Cells.Item(32, 12).Value = 18530800  # L32: 22236820 -> 18530800
$ws.Cells.Item(32, 13).Value = -3281.7463  # M32: -3414.5625 -> -3281.7463
$ws.Cells.Item(32, 14).Value = -18531374  # N32: -22237394 -> -18531374

# Sheet ARM (index 2), row 61
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(61, 8).Value = 4555.9155  # H61: 3898.7585 -> 4555.9155
$ws.Cells.Item(61, 9).Value = 4631.391  # I61: 3134.4614 -> 4631.391
$ws.Cells.Item(61, 11).Value = 4631.391  # K61: 3134.4614 -> 4631.391
$ws.Cells.Item(61, 13).Value = -4419.391  # M61: -2922.4614 -> -4419.391

# Sheet ARM (index 2), row 74
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(74, 8).Value = 4628.9  # H74: 4006.0417 -> 4628.9
$ws.Cells.Item(74, 9).Value = 4643.533  # I74: 3853.6843 -> 4643.533
$ws.Cells.Item(74, 11).Value = 4643.533  # K74: 3853.6843 -> 4643.533
$ws.Cells.Item(74, 13).Value = -3769.533  # M74: -2979.6843 -> -3769.533

# Sheet ARM (index 2), row 77
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(77, 8).Value = 4628.9  # H77: 4006.0417 -> 4628.9
$ws.Cells.Item(77, 9).Value = 4643.533  # I77: 3853.6843 -> 4643.533
$ws.Cells.Item(77, 11).Value = 23217.665  # K77: 19268.4215 -> 23217.665
$ws.Cells.Item(77, 13).Value = -18849.665  # M77: -14900.4215 -> -18849.665

# Sheet ARM (index 2), row 110
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(110, 8).Value = 5847  # H110: 5056.9355 -> 5847
$ws.Cells.Item(110, 9).Value = 2589.6667  # I110: 2003.5714 -> 2589.6667
$ws.Cells.Item(110, 11).Value = 2589.6667  # K110: 2003.5714 -> 2589.6667
$ws.Cells.Item(110, 13).Value = -544.6667000000002  # M110: 41.42859999999996 -> -544.6667000000002

# Sheet ARM (index 2), row 122
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(122, 8).Value = 2270.7856  # H122: 2759.4167 -> 2270.7856
$ws.Cells.Item(122, 9).Value = 1907  # I122: 2373.9092 -> 1907
$ws.Cells.Item(122, 11).Value = 5721  # K122: 7121.7276 -> 5721
$ws.Cells.Item(122, 13).Value = -3271  # M122: -4671.7276 -> -3271

# Sheet ARM (index 2), row 132
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(132, 8).Value = 1905499  # H132: 942164.6 -> 1905499
$ws.Cells.Item(132, 9).Value = 2319639.2  # I132: 1015135.9 -> 2319639.2
$ws.Cells.Item(132, 11).Value = 6958917.600000001  # K132: 3045407.7 -> 6958917.600000001
$ws.Cells.Item(132, 13).Value = -6956387.600000001  # M132: -3042877.7 -> -6956387.600000001

# Sheet ARM (index 2), row 136
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(136, 8).Value = 4555.9155  # H136: 3898.7585 -> 4555.9155
$ws.Cells.Item(136, 9).Value = 4631.391  # I136: 3134.4614 -> 4631.391
$ws.Cells.Item(136, 11).Value = 13894.173  # K136: 9403.3842 -> 13894.173
$ws.Cells.Item(136, 13).Value = -11344.173  # M136: -6853.3842 -> -11344.173

# Sheet BSM (index 3), row 99
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(99, 8).Value = 10081.683  # H99: 10067.049 -> 10081.683
$ws.Cells.Item(99, 9).Value = 9914.870999999999  # I99: 9898.718999999999 -> 9914.870999999999
$ws.Cells.Item(99, 10).Value = 10598.8  # J99: 10665.556 -> 10598.8
$ws.Cells.Item(99, 11).Value = 9914.870999999999  # K99: 9898.718999999999 -> 9914.870999999999
$ws.Cells.Item(99, 12).Value = 10598.8  # L99: 10665.556 -> 10598.8
$ws.Cells.Item(99, 13).Value = -8416.870999999999  # M99: -8400.718999999999 -> -8416.870999999999
$ws.Cells.Item(99, 14).Value = -13594.8  # N99: -13661.556 -> -13594.8

# Sheet BSM (index 3), row 107
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(107, 8).Value = 1965011.4  # H107: 2004348 -> 1965011.4
$ws.Cells.Item(107, 9).Value = 2384466.5  # I107: 2503635.2 -> 2384466.5
$ws.Cells.Item(107, 10).Value = 7554.5557  # J107: 7198.9 -> 7554.5557
$ws.Cells.Item(107, 11).Value = 2384466.5  # K107: 2503635.2 -> 2384466.5
$ws.Cells.Item(107, 12).Value = 7554.5557  # L107: 7198.9 -> 7554.5557
$ws.Cells.Item(107, 13).Value = -2382546.5  # M107: -2501715.2 -> -2382546.5
$ws.Cells.Item(107, 14).Value = -11394.5557  # N107: -11038.9 -> -11394.5557

# Sheet BSM (index 3), row 134
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(134, 8).Value = 7072.9736  # H134: 6356.5347 -> 7072.9736
$ws.Cells.Item(134, 9).Value = 6206.033  # I134: 5449.6855 -> 6206.033
$ws.Cells.Item(134, 11).Value = 18618.099  # K134: 16349.0565 -> 18618.099
$ws.Cells.Item(134, 13).Value = -16083.099  # M134: -13814.0565 -> -16083.099

# Sheet CRP (index 4), row 58
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(58, 8).Value = 47629124  # H58: 47630092 -> 47629124
$ws.Cells.Item(58, 10).Value = 12178.182  # J58: 14021.637 -> 12178.182
$ws.Cells.Item(58, 12).Value = 12178.182  # L58: 14021.637 -> 12178.182
$ws.Cells.Item(58, 14).Value = -12584.182  # N58: -14427.637 -> -12584.182

# Sheet CRP (index 4), row 59
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(59, 8).Value = 33455.668  # H59: 33888.375 -> 33455.668
$ws.Cells.Item(59, 9).Value = 29997  # I59: 30000 -> 29997
$ws.Cells.Item(59, 11).Value = 29997  # K59: 30000 -> 29997
$ws.Cells.Item(59, 13).Value = -28852  # M59: -28855 -> -28852

# Sheet CRP (index 4), row 132
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(132, 8).Value = 9593.700000000001  # H132: 8371.087 -> 9593.700000000001
$ws.Cells.Item(132, 9).Value = 5223  # I132: 4418.3335 -> 5223
$ws.Cells.Item(132, 10).Value = 34361  # J132: 49875 -> 34361
$ws.Cells.Item(132, 11).Value = 15669  # K132: 13255.0005 -> 15669
$ws.Cells.Item(132, 12).Value = 103083  # L132: 149625 -> 103083
$ws.Cells.Item(132, 13).Value = -13139  # M132: -10725.0005 -> -13139
$ws.Cells.Item(132, 14).Value = -108143  # N132: -154685 -> -108143

# Sheet CRP (index 4), row 134
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(134, 8).Value = 95251400  # H134: 127000880 -> 95251400
$ws.Cells.Item(134, 9).Value = 142865380  # I134: 190486530 -> 142865380
$ws.Cells.Item(134, 10).Value = 23437.5  # J134: 29583.334 -> 23437.5
$ws.Cells.Item(134, 11).Value = 428596140  # K134: 571459590 -> 428596140
$ws.Cells.Item(134, 12).Value = 70312.5  # L134: 88750.00199999999 -> 70312.5
$ws.Cells.Item(134, 13).Value = -428593605  # M134: -571457055 -> -428593605
$ws.Cells.Item(134, 14).Value = -75382.5  # N134: -93820.00199999999 -> -75382.5

# Sheet CRP (index 4), row 136
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(136, 8).Value = 47629124  # H136: 47630092 -> 47629124
$ws.Cells.Item(136, 10).Value = 12178.182  # J136: 14021.637 -> 12178.182
$ws.Cells.Item(136, 12).Value = 36534.546  # L136: 42064.911 -> 36534.546
$ws.Cells.Item(136, 14).Value = -41634.546  # N136: -47164.911 -> -41634.546

# Sheet CUL (index 5), row 109
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(109, 8).Value = 8636.362999999999  # H109: 7853.6816 -> 8636.362999999999
$ws.Cells.Item(109, 9).Value = 5312.625  # I109: 4575.4707 -> 5312.625
$ws.Cells.Item(109, 10).Value = 17499.666  # J109: 18999.6 -> 17499.666
$ws.Cells.Item(109, 11).Value = 15937.875  # K109: 13726.4121 -> 15937.875
$ws.Cells.Item(109, 12).Value = 52498.99800000001  # L109: 56998.8 -> 52498.99800000001
$ws.Cells.Item(109, 13).Value = -14897.875  # M109: -12686.4121 -> -14897.875
$ws.Cells.Item(109, 14).Value = -54578.99800000001  # N109: -59078.8 -> -54578.99800000001

# Sheet CUL (index 5), row 132
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(132, 8).Value = 3003.3333  # H132: 3188.2307 -> 3003.3333
$ws.Cells.Item(132, 9).Value = 1703.5  # I132: 1704 -> 1703.5
$ws.Cells.Item(132, 10).Value = 3203.3076  # J132: 3311.9167 -> 3203.3076
$ws.Cells.Item(132, 11).Value = 15331.5  # K132: 15336 -> 15331.5
$ws.Cells.Item(132, 12).Value = 28829.7684  # L132: 29807.2503 -> 28829.7684
$ws.Cells.Item(132, 13).Value = -12801.5  # M132: -12806 -> -12801.5
$ws.Cells.Item(132, 14).Value = -33889.7684  # N132: -34867.2503 -> -33889.7684

# Sheet GSM (index 6), row 113
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(113, 8).Value = 12562.25  # H113: 14285.571 -> 12562.25
$ws.Cells.Item(113, 9).Value = 5250  # I113: 10000 -> 5250
$ws.Cells.Item(113, 10).Value = 14999.667  # J113: 14999.833 -> 14999.667
$ws.Cells.Item(113, 11).Value = 5250  # K113: 10000 -> 5250
$ws.Cells.Item(113, 12).Value = 14999.667  # L113: 14999.833 -> 14999.667
$ws.Cells.Item(113, 13).Value = -3080  # M113: -7830 -> -3080
$ws.Cells.Item(113, 14).Value = -19339.667  # N113: -19339.833 -> -19339.667

# Sheet GSM (index 6), row 122
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(122, 8).Value = 3919.0732  # H122: 4101.3413 -> 3919.0732
$ws.Cells.Item(122, 9).Value = 2827.0857  # I122: 2882.7646 -> 2827.0857
$ws.Cells.Item(122, 10).Value = 10289  # J122: 10020.143 -> 10289
$ws.Cells.Item(122, 11).Value = 8481.257100000001  # K122: 8648.293799999999 -> 8481.257100000001
$ws.Cells.Item(122, 12).Value = 30867  # L122: 30060.429 -> 30867
$ws.Cells.Item(122, 13).Value = -6031.257100000001  # M122: -6198.293799999999 -> -6031.257100000001
$ws.Cells.Item(122, 14).Value = -35767  # N122: -34960.429 -> -35767

# Sheet GSM (index 6), row 132
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(132, 8).Value = 4931.68  # H132: 5023.6733 -> 4931.68
$ws.Cells.Item(132, 9).Value = 4827.4326  # I132: 4949.75 -> 4827.4326
$ws.Cells.Item(132, 11).Value = 14482.2978  # K132: 14849.25 -> 14482.2978
$ws.Cells.Item(132, 13).Value = -11952.2978  # M132: -12319.25 -> -11952.2978

# Sheet GSM (index 6), row 136
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(136, 8).Value = 15518.346  # H136: 15419.714 -> 15518.346
$ws.Cells.Item(136, 10).Value = 15518.346  # J136: 15419.714 -> 15518.346
$ws.Cells.Item(136, 12).Value = 46555.038  # L136: 46259.142 -> 46555.038
$ws.Cells.Item(136, 14).Value = -51655.038  # N136: -51359.142 -> -51655.038

# Sheet LTW (index 7), row 61
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(61, 8).Value = 7675.1816  # H61: 7921.7617 -> 7675.1816
$ws.Cells.Item(61, 9).Value = 5764.4375  # I61: 5982.2666 -> 5764.4375
$ws.Cells.Item(61, 11).Value = 5764.4375  # K61: 5982.2666 -> 5764.4375
$ws.Cells.Item(61, 13).Value = -5562.4375  # M61: -5780.2666 -> -5562.4375

# Sheet LTW (index 7), row 113
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(113, 8).Value = 7675.1816  # H113: 7921.7617 -> 7675.1816
$ws.Cells.Item(113, 9).Value = 5764.4375  # I113: 5982.2666 -> 5764.4375
$ws.Cells.Item(113, 11).Value = 5764.4375  # K113: 5982.2666 -> 5764.4375
$ws.Cells.Item(113, 13).Value = -3594.4375  # M113: -3812.2666 -> -3594.4375

# Sheet LTW (index 7), row 132
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(132, 8).Value = 6806279.5  # H132: 7522287 -> 6806279.5
$ws.Cells.Item(132, 9).Value = 9525818  # I132: 10206154 -> 9525818
$ws.Cells.Item(132, 10).Value = 7433.3335  # J132: 7459 -> 7433.3335
$ws.Cells.Item(132, 11).Value = 28577454  # K132: 30618462 -> 28577454
$ws.Cells.Item(132, 12).Value = 22300.0005  # L132: 22377 -> 22300.0005
$ws.Cells.Item(132, 13).Value = -28574924  # M132: -30615932 -> -28574924
$ws.Cells.Item(132, 14).Value = -27360.0005  # N132: -27437 -> -27360.0005

# Sheet LTW (index 7), row 139
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(139, 8).Value = 0  # H139: 149000 -> 0
$ws.Cells.Item(139, 10).Value = 0  # J139: 149000 -> 0
$ws.Cells.Item(139, 12).Value = 0  # L139: 149000 -> 0
$ws.Cells.Item(139, 14).ClearContents()  # N139: -159280 -> (removed)

# Sheet WVR (index 8), row 107
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(107, 8).Value = 940.4  # H107: 804.6087 -> 940.4
$ws.Cells.Item(107, 9).Value = 1023.1818  # I107: 805.53845 -> 1023.1818
$ws.Cells.Item(107, 10).Value = 839.2222  # J107: 803.4 -> 839.2222
$ws.Cells.Item(107, 11).Value = 3069.5454  # K107: 2416.61535 -> 3069.5454
$ws.Cells.Item(107, 12).Value = 2517.6666  # L107: 2410.2 -> 2517.6666
$ws.Cells.Item(107, 13).Value = -1149.5454  # M107: -496.61535 -> -1149.5454
$ws.Cells.Item(107, 14).Value = -6357.6666  # N107: -6250.2 -> -6357.6666

# Sheet WVR (index 8), row 122
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(122, 8).Value = 5629.091  # H122: 6026.122 -> 5629.091
$ws.Cells.Item(122, 9).Value = 1026.5294  # I122: 1080.5625 -> 1026.5294
$ws.Cells.Item(122, 10).Value = 21277.8  # J122: 23610.334 -> 21277.8
$ws.Cells.Item(122, 11).Value = 3079.5882  # K122: 3241.6875 -> 3079.5882
$ws.Cells.Item(122, 12).Value = 63833.39999999999  # L122: 70831.00199999999 -> 63833.39999999999
$ws.Cells.Item(122, 13).Value = -629.5881999999997  # M122: -791.6875 -> -629.5881999999997
$ws.Cells.Item(122, 14).Value = -68733.39999999999  # N122: -75731.00199999999 -> -68733.39999999999

# Sheet WVR (index 8), row 132
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(132, 8).Value = 16678289  # H132: 17790162 -> 16678289
$ws.Cells.Item(132, 9).Value = 13895073  # I132: 15158071 -> 13895073
$ws.Cells.Item(132, 10).Value = 25027936  # J132: 25028412 -> 25027936
$ws.Cells.Item(132, 11).Value = 41685219  # K132: 45474213 -> 41685219
$ws.Cells.Item(132, 12).Value = 75083808  # L132: 75085236 -> 75083808
$ws.Cells.Item(132, 13).Value = -41682689  # M132: -45471683 -> -41682689
$ws.Cells.Item(132, 14).Value = -75088868  # N132: -75090296 -> -75088868

# Sheet WVR (index 8), row 136
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(136, 8).Value = 18534032  # H136: 15164377 -> 18534032
$ws.Cells.Item(136, 9).Value = 41689696  # I136: 29428278 -> 41689696
$ws.Cells.Item(136, 10).Value = 9501.4  # J136: 8982.25 -> 9501.4
$ws.Cells.Item(136, 11).Value = 125069088  # K136: 88284834 -> 125069088
$ws.Cells.Item(136, 12).Value = 28504.2  # L136: 26946.75 -> 28504.2
$ws.Cells.Item(136, 13).Value = -125066538  # M136: -88282284 -> -125066538
$ws.Cells.Item(136, 14).Value = -33604.2  # N136: -32046.75 -> -33604.2
